$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows -------------------------------------------------------------
# "Deryneia" / "Cyprus" (row 20): inherit the same border/fill/font that the
# rest of the Cyprus rows use (copy the look of the previous row, column B,
# which already carries that plain data style).
$ws.Range("B19").Copy()
$ws.Range("A20").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B19:B19").Copy()
$ws.Range("B20").PasteSpecial(-4122)

$ws.Range("A20").Value = "Deryneia"
$ws.Range("B20").Value = "Cyprus"

# "Polis" / "Cyprus" (row 21)
$ws.Range("B19").Copy()
$ws.Range("B21").PasteSpecial(-4122)

$ws.Range("A21").Value = "Polis"
$ws.Range("B21").Value = "Cyprus"

# --- Alignment --------------------------------------------------------------
# Every data cell (rows 2-21) gets horizontal centering added on top of the
# existing vertical centering / wrap.
$ws.Range("A2:B19").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B20:B21").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A20").HorizontalAlignment = -4108      # xlCenter

# --- Distinct styling for the two new city cells ---------------------------
# "Deryneia" keeps the table look but with an explicit theme font color.
$ws.Range("A20").Font.ThemeColor = 1

# "Polis" drops the table border/fill entirely - plain centered text only.
$ws.Range("A21").HorizontalAlignment = -4108      # xlCenter
$ws.Range("A21").VerticalAlignment = -4108        # xlCenter

$ws.Range("A1:B21").Select()
